$d = $word.ActiveDocument

# --- 1. Merge "Craig Ellis, Andrew Knueven" into a single run, dropping the
#        spell-check proofErr wrapper around "Knueven". The engine leaves a
#        stray trailing proofErr behind when an edit's end boundary lands on
#        the very last position of a paragraph, so we pad with a throwaway
#        character, do the merge, then trim the pad separately.
$p1 = $d.Paragraphs(1)
$padRange = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$padRange.InsertAfter("Z")

$mergeRange = $d.Content
$mergeRange.Find.Execute("Craig Ellis, Andrew KnuevenZ", $true, $false, $false, $false, $false, $true, 1, $false, "Craig Ellis, Andrew KnuevenZ", 2) | Out-Null

$p1b = $d.Paragraphs(1)
$trimRange = $d.Range($p1b.Range.End - 2, $p1b.Range.End - 1)
$trimRange.Delete()

# --- 2. Move the "_GoBack" bookmark from the end of the "Compare and
#        Contrast" section to its new spot inside the Decision Tree
#        paragraph, right after "...which are the IDs of the ".
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$findRange = $d.Content
$findRange.Find.Execute("the IDs of the ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkPos = $findRange.End
$newBookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange) | Out-Null

# --- 3. Append the new "Results" section at the end of the document.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newLastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertRange = $newLastPara.Range
$insertRange.Collapse(0)
$insertRange.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t>Results</w:t></w:r></w:p><w:p><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:tab/></w:r><w:r><w:t>Our decision trees and RAP files were tuned using an iterative process. The “data” folder contains 5 decision trees and 7 RAPs, each tuning different bits of the architecture and adding preconditions until we reached what we believe is the most efficient version of each system. To test the per</w:t></w:r><w:r><w:t xml:space="preserve">formance of each, we hooked up our controllers </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DTPacMan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RAPPacMan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to the Executor class and created a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>runExperiment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method that runs the game a set number of times in asynchronous mode with the visuals off, so it doesn’t take an inordinate amount of time to run a large amount of trials. We selected 100 runs as a good estimate of the average score obtained by each system – the experiment collects and prints the score for each run and then the average at the end. Results for the final iteration of each system can be found in “results.txt” of the root directory. The DT and RAP performed about the same – since the decision tree and RAP are both </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>really only</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> giving Pac-Man one move to do at any given point, they perform about the same given the same set of actions to perform. The average score for each of these was about 11,000 over the 100 runs.</w:t></w:r></w:p><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
